$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

# Row 7 (Q6, A col value = 11)
$ws.Range("B7").Value = 0.05422975282401672
$ws.Range("C7").Value = 2.682531801186378
$ws.Range("D7").Value = 18.95522674257797
$ws.Range("E7").Value = 4.353760069477643
$ws.Range("F7").Value = 4.41186013140652
$ws.Range("G7").Value = 38

# Row 8 (Q7, A col value = 12)
$ws.Range("B8").Value = 0.09175281835627994
$ws.Range("C8").Value = 2.663306588085196
$ws.Range("D8").Value = 17.83211946726106
$ws.Range("E8").Value = 4.222809428243365
$ws.Range("F8").Value = 4.280047159788971
$ws.Range("G8").Value = 37

# Row 9 (Q8, A col value = 13)
$ws.Range("B9").Value = 0.07118176867328285
$ws.Range("C9").Value = 3.846179246344002
$ws.Range("D9").Value = 30.36400025844973
$ws.Range("E9").Value = 5.51035391408299
$ws.Range("F9").Value = 5.653032108082237
$ws.Range("G9").Value = 20

# Row 10 (Q9, A col value = 14)
$ws.Range("B10").Value = 0.0002549338085925304
$ws.Range("C10").Value = 4.420675973149636
$ws.Range("D10").Value = 37.99070145917987
$ws.Range("E10").Value = 6.163659745571609
$ws.Range("F10").Value = 6.415340456829838
$ws.Range("G10").Value = 13

# Row 11 (A col value = 15)
$ws.Range("B11").Value = -0.1404721957463057
$ws.Range("C11").Value = 3.824781345671012
$ws.Range("D11").Value = 24.88086598134253
$ws.Range("E11").Value = 4.988072371301617
$ws.Range("F11").Value = 5.574622581794747
$ws.Range("G11").Value = 5
